$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the "%%b" -> "%%B" typo in the Buttons example output text (old row 198, col B).
$ws.Range("B198").Value2 = "Toto jsou vaše možnosti: %%B teď byl převezen=teď byl převezen do nemocnice;brzy propustí= pacienta brzy propustí z nemocnice "

# 2) Insert 5 new rows (200-204) for the new "Foldables" tutorial section,
#    pushing the existing "Kde a jak pokracovat v konverzaci" section (old row 200+) down to row 205+.
$ws.Rows("200:204").Insert()

# Row 200: bold heading "// Foldables"
$ws.Range("A200").Value2 = "// Foldables"
$ws.Range("A200").Characters(1, 3).Font.Bold = $false
$ws.Range("A200").Characters(4, 9).Font.Bold = $true

# Row 201: explanatory comment about foldables
$ws.Range("A201").Value2 = "// Foldables  (rozbalovací texty) se na yačátku yobrayí ve své krátké variantě (npříklad jen jako titulek) a pokud na ně uživatel klikne, krátká forma se nahradí  dlouhou"

# Row 202: example row - button/variable name in col A, output spec with %%F in col B
$ws.Range("A202").Value2 = "foldable"
$ws.Range("B202").Value2 = "Příklad %%Ftitulek=dlouhá forma textu; druhý titulek= a velmi dlouhá forma druhého titulku"

# Row 203: the (misspelled, as in source) label "roybalovací text"
$ws.Range("A203").Value2 = "roybalovací text"

# Row 204 stays blank (spacer row), matching the surrounding tutorial formatting.

# Apply the same column-A "comment" style used by the neighboring rows (195/200/201/203/205...)
# so the new rows visually match the rest of the section.
$ws.Range("A200:A204").Style = $ws.Range("A199").Style

# 3) Move the selection to reflect where editing ended up.
$ws.Range("B198").Select()
